# Applies the "Updated cryptos list" data refresh described by the commit.
# Most cells are plain text (inlineStr) so setting .Value directly is fine.
# Cells in column D sometimes hold a numeric-looking piece of text (e.g. "1.01")
# which real Excel would silently convert into a *number* the moment it is
# assigned to a normally-formatted cell. To keep those values as literal TEXT
# (matching the source file), we briefly mark the cell as Text ("@") before
# writing the value, then clear the formatting again so no stray number format
# is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "26.052.29"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.638.98"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.55%  "
Set-TextValue "D5" "214.75"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("E8").Value = "  -1.63%  "
Set-TextValue "D9" "0.0628"
$ws.Range("E9").Value = "  -1.29%  "
Set-TextValue "D10" "18.71"
$ws.Range("E10").Value = "  -4.52%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D12" "4.21"
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.642.54"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("E14").Value = "  -2.10%  "
Set-TextValue "D15" "62.33"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("D17").Value = "26.069.81"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("E18").Value = "  +0.48%  "
Set-TextValue "D19" "191.43"
$ws.Range("E19").Value = "  -0.68%  "
Set-TextValue "D20" "4.28"
$ws.Range("E20").Value = "  -1.48%  "
Set-TextValue "D21" "9.64"
$ws.Range("E21").Value = "  -2.66%  "
Set-TextValue "D22" "6.16"
$ws.Range("E22").Value = "  -1.56%  "
Set-TextValue "D23" "0.133"
$ws.Range("E23").Value = "  +2.09%  "
Set-TextValue "D24" "144.09"
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("E26").Value = "  -0.98%  "
Set-TextValue "D27" "6.79"
$ws.Range("E27").Value = "  -1.27%  "
Set-TextValue "D28" "15.27"
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("E29").Value = "  -0.21%  "
Set-TextValue "D30" "0.0488"
$ws.Range("E30").Value = "  -2.60%  "
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("E32").Value = "  -3.12%  "
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("E34").Value = "  -0.83%  "
Set-TextValue "D35" "0.879"
$ws.Range("E35").Value = "  -2.39%  "
$ws.Range("D36").Value = "1.129.51"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("E37").Value = "  -0.22%  "
Set-TextValue "D38" "0.527"
$ws.Range("E38").Value = "  -2.92%  "
$ws.Range("E39").Value = "  -1.03%  "
Set-TextValue "D40" "99.01"
$ws.Range("E40").Value = "  -0.27%  "
Set-TextValue "D41" "0.786"
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("E42").Value = "  -2.95%  "
$ws.Range("E43").Value = "  -0.56%  "
Set-TextValue "D44" "55.55"
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("E45").Value = "  -0.21%  "
Set-TextValue "D46" "1.49"
$ws.Range("E46").Value = "  +1.33%  "
Set-TextValue "D47" "0.414"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("E48").Value = "  -1.09%  "
Set-TextValue "D49" "1.01"
Set-TextValue "D50" "0.0931"
$ws.Range("E50").Value = "  -2.88%  "
$ws.Range("E51").Value = "  -0.60%  "
